$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style used by unstyled data cells (e.g. A2), used to reset styling after
# forcing NumberFormat="@" on Price (D) column cells -- avoids Excel silently
# re-typing numeric-looking price strings (e.g. "591.80") as floats, which would
# both change the stored type and lose significant trailing zeros.
$normalStyle = $ws.Range("A2").Style

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.378.51'
$ws.Range('D2').Style = $normalStyle
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.508.15'
$ws.Range('D3').Style = $normalStyle
$ws.Range('E3').Value = '  +0.78%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '591.80'
$ws.Range('D5').Style = $normalStyle
$ws.Range('E5').Value = '  +1.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.64'
$ws.Range('D6').Style = $normalStyle
$ws.Range('E6').Value = '  +0.47%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  +0.56%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.60'
$ws.Range('D9').Style = $normalStyle
$ws.Range('E9').Value = '  +6.07%  '
$ws.Range('E10').Value = '  +1.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.388'
$ws.Range('D11').Style = $normalStyle
$ws.Range('E11').Value = '  +3.66%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.106.71'
$ws.Range('D12').Style = $normalStyle
$ws.Range('E12').Value = '  +0.64%  '
$ws.Range('E13').Value = '  +1.05%  '
$ws.Range('E14').Value = '  +1.08%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.508.39'
$ws.Range('D15').Style = $normalStyle
$ws.Range('E15').Value = '  +0.65%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '25.80'
$ws.Range('D16').Style = $normalStyle
$ws.Range('E16').Value = '  +3.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.357.90'
$ws.Range('D17').Style = $normalStyle
$ws.Range('E17').Value = '  +0.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.02'
$ws.Range('D18').Style = $normalStyle
$ws.Range('E18').Value = '  +0.77%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.63'
$ws.Range('D19').Style = $normalStyle
$ws.Range('E19').Value = '  -0.35%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.76'
$ws.Range('D20').Style = $normalStyle
$ws.Range('E20').Value = '  +2.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '391.45'
$ws.Range('D21').Style = $normalStyle
$ws.Range('E21').Value = '  +1.51%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.581'
$ws.Range('D22').Style = $normalStyle
$ws.Range('E22').Value = '  +3.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.648.08'
$ws.Range('D23').Style = $normalStyle
$ws.Range('E23').Value = '  +0.65%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '74.50'
$ws.Range('D24').Style = $normalStyle
$ws.Range('E24').Value = '  +0.72%  '
$ws.Range('E25').Value = '  +0.12%  '
$ws.Range('E26').Value = '  -0.55%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000117'
$ws.Range('D27').Style = $normalStyle
$ws.Range('E27').Value = '  +4.07%  '
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.43'
$ws.Range('D29').Style = $normalStyle
$ws.Range('E29').Value = '  +0.72%  '
$ws.Range('E30').Value = '  +2.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.23'
$ws.Range('D31').Style = $normalStyle
$ws.Range('E31').Value = '  +0.72%  '
$ws.Range('E32').Value = '  -4.16%  '
$ws.Range('E33').Value = '  +7.71%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.534.40'
$ws.Range('D34').Style = $normalStyle
$ws.Range('E34').Value = '  +0.98%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '23.37'
$ws.Range('D36').Style = $normalStyle
$ws.Range('E36').Value = '  +0.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.34'
$ws.Range('D37').Style = $normalStyle
$ws.Range('E37').Value = '  +0.94%  '
$ws.Range('E38').Value = '  +1.89%  '
$ws.Range('E39').Value = '  +2.09%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '166.68'
$ws.Range('D40').Style = $normalStyle
$ws.Range('E40').Value = '  +2.54%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0788'
$ws.Range('D41').Style = $normalStyle
$ws.Range('E41').Value = '  +1.53%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.811'
$ws.Range('D42').Style = $normalStyle
$ws.Range('E42').Value = '  +1.28%  '
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '25.17'
$ws.Range('D44').Style = $normalStyle
$ws.Range('E44').Value = '  -0.76%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.45'
$ws.Range('D45').Style = $normalStyle
$ws.Range('E45').Value = '  +1.67%  '
$ws.Range('E46').Value = '  +0.23%  '
$ws.Range('E47').Value = '  -1.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.80'
$ws.Range('D48').Style = $normalStyle
$ws.Range('E48').Value = '  +1.02%  '
$ws.Range('B49').Value = 'SuiNetwork'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.911'
$ws.Range('D49').Style = $normalStyle
$ws.Range('E49').Value = '  +1.37%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.374.27'
$ws.Range('D50').Style = $normalStyle
$ws.Range('E50').Value = '  -3.77%  '
$ws.Range('E51').Value = '  +0.52%  '
